# Scheduled runner update: refresh cached Market Board price/profit figures
# (currentAveragePrice / NQ / HQ / Leve*Price / Leve*Profit columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW leve-profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1035.625
$ws.Range("I28").Value = 749.36365
$ws.Range("J28").Value = 1665.4
$ws.Range("K28").Value = 749.36365
$ws.Range("L28").Value = 1665.4
$ws.Range("M28").Value = -264.36365
$ws.Range("N28").Value = -2635.4

$ws.Range("H32").Value = 6932.6665
$ws.Range("I32").Value = 6801
$ws.Range("J32").Value = 6998.5
$ws.Range("K32").Value = 6801
$ws.Range("L32").Value = 6998.5
$ws.Range("M32").Value = -6475
$ws.Range("N32").Value = -7650.5

$ws.Range("H106").Value = 4029.2144
$ws.Range("I106").Value = 3308.5833
$ws.Range("J106").Value = 4569.6875
$ws.Range("K106").Value = 3308.5833
$ws.Range("L106").Value = 4569.6875
$ws.Range("M106").Value = -2677.5833
$ws.Range("N106").Value = -5831.6875

$ws.Range("H111").Value = 2999
$ws.Range("I111").Value = 2999
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 8997
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -5930

$ws.Range("H132").Value = 66210.625
$ws.Range("I132").Value = 70366
$ws.Range("J132").Value = 3880
$ws.Range("K132").Value = 211098
$ws.Range("L132").Value = 11640
$ws.Range("M132").Value = -208568
$ws.Range("N132").Value = -16700

$ws.Range("H138").Value = 6246.8906
$ws.Range("I138").Value = 2925.125
$ws.Range("J138").Value = 7179.316
$ws.Range("K138").Value = 8775.375
$ws.Range("L138").Value = 21537.948
$ws.Range("M138").Value = -3635.375
$ws.Range("N138").Value = -31817.948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10200.25
$ws.Range("I2").Value = 16024.75
$ws.Range("J2").Value = 4375.75
$ws.Range("K2").Value = 16024.75
$ws.Range("L2").Value = 4375.75
$ws.Range("M2").Value = -15911.75
$ws.Range("N2").Value = -4601.75

$ws.Range("H32").Value = 2344.2927
$ws.Range("I32").Value = 2277.9
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 2277.9
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -1990.9
$ws.Range("N32").Value = -5574

$ws.Range("H45").Value = 1567.8
$ws.Range("I45").Value = 1570.6923
$ws.Range("J45").Value = 1549
$ws.Range("K45").Value = 1570.6923
$ws.Range("L45").Value = 1549
$ws.Range("M45").Value = -1193.6923
$ws.Range("N45").Value = -2303

$ws.Range("H61").Value = 2565.0908
$ws.Range("I61").Value = 2524.111
$ws.Range("J61").Value = 2749.5
$ws.Range("K61").Value = 2524.111
$ws.Range("L61").Value = 2749.5
$ws.Range("M61").Value = -2312.111
$ws.Range("N61").Value = -3173.5

$ws.Range("H74").Value = 753965.25
$ws.Range("I74").Value = 3338.257
$ws.Range("J74").Value = 13889938
$ws.Range("K74").Value = 3338.257
$ws.Range("L74").Value = 13889938
$ws.Range("M74").Value = -2464.257
$ws.Range("N74").Value = -13891686

$ws.Range("H77").Value = 753965.25
$ws.Range("I77").Value = 3338.257
$ws.Range("J77").Value = 13889938
$ws.Range("K77").Value = 16691.285
$ws.Range("L77").Value = 69449690
$ws.Range("M77").Value = -12323.285
$ws.Range("N77").Value = -69458426

$ws.Range("H116").Value = 10200.25
$ws.Range("I116").Value = 16024.75
$ws.Range("J116").Value = 4375.75
$ws.Range("K116").Value = 16024.75
$ws.Range("L116").Value = 4375.75
$ws.Range("M116").Value = -13730.75
$ws.Range("N116").Value = -8963.75

$ws.Range("H136").Value = 2565.0908
$ws.Range("I136").Value = 2524.111
$ws.Range("J136").Value = 2749.5
$ws.Range("K136").Value = 7572.333
$ws.Range("L136").Value = 8248.5
$ws.Range("M136").Value = -5022.333
$ws.Range("N136").Value = -13348.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10200.25
$ws.Range("I3").Value = 16024.75
$ws.Range("J3").Value = 4375.75
$ws.Range("K3").Value = 16024.75
$ws.Range("L3").Value = 4375.75
$ws.Range("M3").Value = -15910.75
$ws.Range("N3").Value = -4603.75

$ws.Range("H134").Value = 27779594
$ws.Range("I134").Value = 41667984
$ws.Range("J134").Value = 2814
$ws.Range("K134").Value = 125003952
$ws.Range("L134").Value = 8442
$ws.Range("M134").Value = -125001417
$ws.Range("N134").Value = -13512

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2637
$ws.Range("I16").Value = 2637
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2637
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2350

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0

$ws.Range("H113").Value = 2637
$ws.Range("I113").Value = 2637
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2637
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -467

$ws.Range("H134").Value = 2732.5
$ws.Range("I134").Value = 2999
$ws.Range("J134").Value = 2599.25
$ws.Range("K134").Value = 8997
$ws.Range("L134").Value = 7797.75
$ws.Range("M134").Value = -6462
$ws.Range("N134").Value = -12867.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1024.6364
$ws.Range("I18").Value = 346.875
$ws.Range("J18").Value = 2832
$ws.Range("K18").Value = 1040.625
$ws.Range("L18").Value = 8496
$ws.Range("M18").Value = -871.625
$ws.Range("N18").Value = -8834

$ws.Range("H34").Value = 5104.5654
$ws.Range("I34").Value = 745.6667
$ws.Range("J34").Value = 6643
$ws.Range("K34").Value = 2237.0001
$ws.Range("L34").Value = 19929
$ws.Range("M34").Value = -2153.0001
$ws.Range("N34").Value = -20097

$ws.Range("H39").Value = 6269.375
$ws.Range("I39").Value = 1519.6
$ws.Range("J39").Value = 8428.362999999999
$ws.Range("K39").Value = 4558.799999999999
$ws.Range("L39").Value = 25285.089
$ws.Range("M39").Value = -4264.799999999999
$ws.Range("N39").Value = -25873.089

$ws.Range("H55").Value = 5355.8887
$ws.Range("I55").Value = 99
$ws.Range("J55").Value = 6857.857
$ws.Range("K55").Value = 297
$ws.Range("L55").Value = 20573.571
$ws.Range("M55").Value = -120
$ws.Range("N55").Value = -20927.571

$ws.Range("H131").Value = 48877.855
$ws.Range("I131").Value = 77670.92
$ws.Range("J131").Value = 35970.62
$ws.Range("K131").Value = 233012.76
$ws.Range("L131").Value = 107911.86
$ws.Range("M131").Value = -227972.76
$ws.Range("N131").Value = -117991.86

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12470
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4479.6
$ws.Range("I16").Value = 4524.75
$ws.Range("J16").Value = 4299
$ws.Range("K16").Value = 4524.75
$ws.Range("L16").Value = 4299
$ws.Range("M16").Value = -4354.75
$ws.Range("N16").Value = -4639

$ws.Range("H46").Value = 2616.8333
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 3040.2
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 3040.2
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -3416.2

$ws.Range("H68").Value = 6105.8
$ws.Range("I68").Value = 4578.3335
$ws.Range("J68").Value = 14125
$ws.Range("K68").Value = 4578.3335
$ws.Range("L68").Value = 14125
$ws.Range("M68").Value = -3829.3335
$ws.Range("N68").Value = -15623

$ws.Range("H71").Value = 6105.8
$ws.Range("I71").Value = 4578.3335
$ws.Range("J71").Value = 14125
$ws.Range("K71").Value = 22891.6675
$ws.Range("L71").Value = 70625
$ws.Range("M71").Value = -19147.6675
$ws.Range("N71").Value = -78113

$ws.Range("H100").Value = 1469.4
$ws.Range("I100").Value = 1649.3334
$ws.Range("J100").Value = 1199.5
$ws.Range("K100").Value = 1649.3334
$ws.Range("L100").Value = 1199.5
$ws.Range("M100").Value = -1108.3334
$ws.Range("N100").Value = -2281.5

$ws.Range("H136").Value = 43482420
$ws.Range("I136").Value = 4096.2354
$ws.Range("J136").Value = 166670990
$ws.Range("K136").Value = 12288.7062
$ws.Range("L136").Value = 500012970
$ws.Range("M136").Value = -9738.706199999999
$ws.Range("N136").Value = -500018070
